$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.678.71"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "'2.900.73"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'586.48"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "'146.62"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'2.899.47"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").Value = "'6.81"
$ws.Range("E10").Value = "  -5.50%  "
$ws.Range("D11").Value = "'0.152"
$ws.Range("E11").Value = "  +5.29%  "
$ws.Range("D12").Value = "'0.433"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "'32.66"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "'3.378.12"
$ws.Range("E16").Value = "  -2.26%  "
$ws.Range("D17").Value = "'61.643.89"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "'6.65"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'2.843.71"
$ws.Range("E19").Value = "  -4.97%  "
$ws.Range("D20").Value = "'437.53"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "'13.31"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "'0.658"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "'6.95"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").Value = "'80.01"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "'11.96"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("D26").Value = "'10.18"
$ws.Range("E26").Value = "  -9.87%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'2.06"
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").Value = "'0.0000111"
$ws.Range("E29").Value = "  +19.08%  "
$ws.Range("D30").Value = "'7.14"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "'2.55"
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").Value = "'0.108"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'25.81"
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").Value = "'0.969"
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").Value = "'3.09"
$ws.Range("E37").Value = "  +3.48%  "
$ws.Range("D38").Value = "'5.49"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").Value = "'49.14"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'1.99"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").Value = "'8.34"
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("D42").Value = "'0.116"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").Value = "'0.270"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("D44").Value = "'38.09"
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'134.95"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'2.683.16"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D47").Value = "'0.0335"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'340.82"
$ws.Range("E49").Value = "  -6.68%  "
$ws.Range("D50").Value = "'0.103"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "'21.98"
$ws.Range("E51").Value = "  -4.66%  "
